$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 785
$ws1.Range("F3").Value = 61
$ws1.Range("F5").Value = 142
$ws1.Range("F7").Value = 162
$ws1.Range("F9").Value = 460
$ws1.Range("F10").Value = 516
$ws1.Range("F12").Value = 11840
$ws1.Range("F13").Value = 5421

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 785
$ws4.Range("F3").Value = 61
$ws4.Range("F7").Value = 142
$ws4.Range("F9").Value = 162
$ws4.Range("F11").Value = 460
$ws4.Range("F12").Value = 516
$ws4.Range("F14").Value = 11840
$ws4.Range("F16").Value = 5421
